$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 / J1, matching the formatting already used by H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I column is 1 for every row except row 8 (value 3)
# J column mirrors column H for every row except row 8 (value 8)
$iValues = @(1,1,1,1,1,1,3,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$jValues = @(5,5,6,6,7,5,8,5,7,5,6,5,4,5,5,6,7,6,6,7,7,5,5,6,6,6,5,5,6,4,5,6,5,5,3,3,2)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
